$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet currently has two "header" rows (1 and 2) that are remnants of a
# larger merged table (truncated "...mation" / "...pompes)" labels, plus
# "Hiver"/"Eté"/"Année" season labels in row 1 and unit labels in row 2).
# Replace them with a single proper header row: idx, idx2, Name, Date Start,
# Date End, (m3/s), (MW1), (MW2), (GWh) Winter, (GWh) Summer, (GWh) Year.

# Remove the old two header rows (data rows shift up to start at row 1)...
$ws.Range("A1:A2").EntireRow.Delete()
# ...then insert one fresh blank row back at the top for the new header.
$ws.Range("A1").EntireRow.Insert()

$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# F1:K1 pick up the 9pt-Arial header styling already used elsewhere in the sheet.
$ws.Range("F1:K1").Font.Size = 9

# Match the author's final selection (row 2, the first data row).
[void]$ws.Range("A2:K2").Select()
